$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.125.62"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "3.133.75"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'581.36"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "'174.13"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'6.45"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").Value = "'0.156"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "'0.482"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").Value = "'37.65"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "'0.123"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").Value = "67.029.64"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "'7.15"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").Value = "3.132.64"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").Value = "'16.42"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "'492.93"
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").Value = "'7.90"
$ws.Range("E21").Value = "  +5.13%  "
$ws.Range("D22").Value = "'84.23"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").Value = "'13.28"
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("D24").Value = "'2.31"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").Value = "'10.38"
$ws.Range("E25").Value = "  +3.73%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'7.97"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").Value = "'28.81"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "'0.115"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").Value = "0.0₃0956"
$ws.Range("E32").Value = "  -5.10%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'5.92"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").Value = "'0.980"
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("D36").Value = "'46.96"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "'2.07"
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'50.12"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "'0.313"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").Value = "'8.58"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "'387.28"
$ws.Range("E42").Value = "  +2.27%  "
$ws.Range("D43").Value = "2.832.97"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "'2.62"
$ws.Range("E44").Value = "  -6.17%  "
$ws.Range("D45").Value = "'0.0355"
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("D46").Value = "'136.07"
$ws.Range("E46").Value = "  +0.66%  "
$ws.Range("D48").Value = "'25.14"
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").Value = "'6.80"
$ws.Range("E51").Value = "  -0.17%  "
